# Switch the "period" column on sheet M from text labels ("1998M01" ...)
# to real date values (first of each month, Jan-Jun 2020), formatted as
# short dates, and fix the downstream numeric formatting that implies.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M")

# Date serials (Excel 1900 date system) for 2020-01-01 .. 2020-06-01.
$dateSerials = @(43831, 43862, 43891, 43922, 43952, 43983)

for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $dateSerials[$i]
}

# Apply a date number format to the first cell, then propagate the same
# style to the rest of the column via a format-only paste so every cell
# in the range shares a single style record.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
